# Generate Report for handback
# The "9c0760a9-f347-488d-8db0-4bdee644a390.md" file has now been handed
# back, so update its Status from "Not yet handed off" to "Handed back"
# on every sheet, and record the handback datetime for each locale.

$wb = $excel.ActiveWorkbook

# Overview sheet: Status columns for zh-cn (B) and de-de (C)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back"
$wsOverview.Range("C3").Value = "Handed back"

# zh-cn sheet: Status (B) and Latest Handback DateTime (G)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Handed back"
$wsZhCn.Range("G3").Value = "2016-01-07 07:47:57"

# de-de sheet: Status (B) and Latest Handback DateTime (G)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Handed back"
$wsDeDe.Range("G3").Value = "2016-01-07 07:48:15"
